$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit: the first "Pinout (pins are labelled on the underside of the
# PCB):" heading (in the very first table cell) gets split by a stray extra
# "m" typed after "underside", with the document's "_GoBack" last-edit
# bookmark left right after that "m" (and consequently removed from its old
# location, the empty paragraph that follows the big table at the end of the
# document).
# ---------------------------------------------------------------------------

$needle     = "Pinout (pins are labelled on the underside of the PCB):"
$firstPart  = "Pinout (pins are labelled on the underside"

$full = $d.Content.Text
$startIdx = $full.IndexOf($needle)
$splitIdx = $startIdx + $firstPart.Length

# Insert the extra "m" right after "...underside" (and before " of the
# PCB):"). This merges into the existing run, which already carries the
# correct bold/underline/font formatting.
$insertRange = $d.Range($splitIdx, $splitIdx)
$insertRange.InsertAfter("m")

# Force the newly-typed "m" to live in its own <w:r> (matching a real
# editing session where the run got split) by round-tripping a formatting
# property that is already explicitly set on it - this changes nothing
# visually but stops the writer from re-merging it with its neighbour.
$mRange = $d.Range($splitIdx, $splitIdx + 1)
$mRange.Font.Bold = 0
$mRange.Font.Bold = 1

# Re-add the "_GoBack" bookmark immediately after the "m" (i.e. right before
# " of the PCB):"). Word only allows one bookmark per name, so adding it here
# moves it away from wherever it used to be (the trailing empty paragraph at
# the very end of the document), cleanly removing it from there.
$bmRange = $d.Range($splitIdx + 1, $splitIdx + 1)
$bmRange.Bookmarks.Add("_GoBack")
